$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Attendance value in B14 (previously "Madison"), keep its format
$ws.Range("B14").Value = $null

# Add a new row of minutes for November 16 2023, matching the formatting
# used by the rest of the table (copy cell-by-cell down from row 14 so we
# don't introduce stray blank cells in columns that should stay empty).
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A14").Copy()
$ws.Range("E15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E16").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

$ws.Range("A15").Value = "November 16 2023"
$ws.Range("B15").Value = "Sedat, David, Madison, Sean, Joseph"
$ws.Range("C15").Value = "1:15PM"
$ws.Range("E15").Value = "Worked on the java files of application"

# Update the view state to match the saved selection
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E20").Select()
